# Actualización automática 2025-11-24 08:30:09
#
# This script reproduces, via Excel COM interop, the insertion of a new
# advisor ("PERDOMO BRIONES JOSÉ ALBERTO") at row 45 of the two detail
# sheets plus a handful of updated sales figures and the resulting
# recomputed summary/compliance numbers on the third sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Updated figures for existing clients/advisors.
$ws1.Range("H21").Value = 495
$ws1.Range("I21").Value = 528.53
$ws1.Range("M21").Value = 2208.95
$ws1.Range("L33").Value = 190.68
$ws1.Range("L43").Value = 517.0599999999999

# Insert a new row for the new advisor, pushing everyone from the old
# row 45 onward down by one row.
$ws1.Rows.Item(45).Insert()

$ws1.Range("A45").Value = "OFICINA-CATAECSA"
$ws1.Range("B45").Value = "PERDOMO BRIONES JOSÉ ALBERTO"
$ws1.Range("C45").Value = 0
$ws1.Range("D45").Value = 0
$ws1.Range("E45").Value = 169.73
$ws1.Range("F45").Value = 0
$ws1.Range("G45").Value = 0
$ws1.Range("H45").Value = 0
$ws1.Range("I45").Value = 0
$ws1.Range("J45").Value = 0
$ws1.Range("K45").Value = 0
$ws1.Range("L45").Value = 413.48
$ws1.Range("M45").Value = 0
$ws1.Range("N45").Value = 0
$ws1.Range("O45").Value = 0
$ws1.Range("P45").Value = 0
$ws1.Range("Q45").Value = 0
$ws1.Range("R45").Value = 0

# The trailing "count of non-zero entries" summary row was pushed from
# row 57 to row 58; refresh its "X de 55" -> "X de 56" labels.
$ws1.Range("C58").Value = "0 de 56"
$ws1.Range("D58").Value = "0 de 56"
$ws1.Range("E58").Value = "2 de 56"
$ws1.Range("F58").Value = "0 de 56"
$ws1.Range("G58").Value = "0 de 56"
$ws1.Range("H58").Value = "1 de 56"
$ws1.Range("I58").Value = "1 de 56"
$ws1.Range("J58").Value = "0 de 56"
$ws1.Range("K58").Value = "0 de 56"
$ws1.Range("L58").Value = "3 de 56"
$ws1.Range("M58").Value = "2 de 56"
$ws1.Range("N58").Value = "0 de 56"
$ws1.Range("O58").Value = "0 de 56"
$ws1.Range("P58").Value = "0 de 56"
$ws1.Range("Q58").Value = "0 de 56"
$ws1.Range("R58").Value = "0 de 56"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F21").Value = 3232.48
$ws2.Range("F33").Value = 216.76
$ws2.Range("F43").Value = 517.0599999999999

# Same row insertion as sheet 1.
$ws2.Rows.Item(45).Insert()

$ws2.Range("A45").Value = "OFICINA-CATAECSA"
$ws2.Range("B45").Value = "PERDOMO BRIONES JOSÉ ALBERTO"
$ws2.Range("C45").Value = 0
$ws2.Range("D45").Value = 0
$ws2.Range("E45").Value = 0
$ws2.Range("F45").Value = 583.21
$ws2.Range("G45").Value = 0

# Totals row moved from row 57 to row 58; only the "noviembre" total
# changes.
$ws2.Range("F58").Value = 5493.93

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 8039.94
$ws3.Range("E3").Value = -8039.94

$ws3.Range("D4").Value = 10987.13
$ws3.Range("E4").Value = 15012.87
$ws3.Range("F4").Value = 0.422581923076923

$ws3.Range("D5").Value = 19027.07
$ws3.Range("E5").Value = 7059.340000000001
$ws3.Range("F5").Value = 0.7293862973095953

# Column D widens from 13 to 14 characters, column F narrows from 24 to
# 23. Excel's ColumnWidth setter adds a fixed ~0.8333 character padding
# on top of the value you assign before it lands in the stored <col>
# width, so back that offset out to land exactly on the target widths.
$ws3.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws3.Columns.Item(6).ColumnWidth = 22.166666666666668
